$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for two new attack-animation columns (Attack2AnimName /
#    Attack3AnimName) between the existing "AttackAnimName" column (L) and
#    "GetHitAnimName"/"DeathAnimName" (old M/N). Move the old M/N data out to
#    O/P first, then reuse M/N for the new columns.
# ---------------------------------------------------------------------------
$oldM = @{}
$oldN = @{}
for ($r = 1; $r -le 12; $r++) {
    $oldM[$r] = $ws.Cells.Item($r, 13).Value2
    $oldN[$r] = $ws.Cells.Item($r, 14).Value2
}

$ws.Range("M1:N12").Clear()

for ($r = 1; $r -le 12; $r++) {
    $ws.Cells.Item($r, 15).Value = $oldM[$r]
    $ws.Cells.Item($r, 16).Value = $oldN[$r]
}

# Rename the old "AttackAnimName" header to "Attack1AnimName" and add headers
# for the two new columns.
$ws.Range("L1").Value = "Attack1AnimName"
$ws.Range("M1").Value = "Attack2AnimName"
$ws.Range("N1").Value = "Attack3AnimName"

# Match the source column width/format for the newly (re)used column M.
$ws.Columns("M").ColumnWidth = 19.375

# ---------------------------------------------------------------------------
# 2. Fill the two new columns for every data row (2-12). Default to "None";
#    Ursacetus gets its real boss combo (set below).
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Value = "None"
    $ws.Cells.Item($r, 14).Value = "None"
}

# ---------------------------------------------------------------------------
# 3. Data updates (balance tweaks + stage-1 boss attack pattern).
# ---------------------------------------------------------------------------

# Onyscidus: MaxHp 70 -> 100
$ws.Range("K2").Value = 100

# Arack: MaxHp 50 -> 30
$ws.Range("K3").Value = 30

# Ceratoferox: MaxHp 300 -> 200
$ws.Range("K6").Value = 200

# Ursacetus (stage 1 boss): MoveSpeed, AttackRange tweaks + new attack combo
$ws.Range("D7").Value = 5
$ws.Range("I7").Value = 5
$ws.Range("L7").Value = "LeftHandAttack"
$ws.Range("M7").Value = "Roar1"
$ws.Range("N7").Value = "LeftFootStompAttack"
$ws.Range("O7").Value = "None"

# ---------------------------------------------------------------------------
# 4. Re-sort the full data range (now A2:P12) ascending by EnemyLevel (B),
#    matching the author re-running the sort after adding the new rows.
# ---------------------------------------------------------------------------
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B12"))
$ws.Sort.SetRange($ws.Range("A1:P12"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# ---------------------------------------------------------------------------
# 5. Misc view-state touch-ups.
# ---------------------------------------------------------------------------
$ws.Range("L16").Select()
